$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): add two new columns for the "Cierre" accounting dynamic ---
$ws.Range("K1").Value = "AccountingSourceCierre"
$ws.Range("L1").Value = "AccountingNameCierre"

# --- Row 2 (data values) ---
$ws.Range("A2").Value = "183"
$ws.Range("C2").Value = """ALCALDIA MUNICIPAL DE IBAGUE PENSIONADOS"""
$ws.Range("E2").Value = "2021"
$ws.Range("I2").Value = """upper('Aplicación de pago por pagaduría') """
$ws.Range("J2").Value = "17/01/2022"
$ws.Range("K2").Value = """'CIERRE'"""
$ws.Range("L2").Value = """upper('Causación fianza cierre de periodo')"""

# --- Update the view so the new columns (K:L) are visible and L1 is selected ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("L1").Select() | Out-Null
